# Updated cryptos list on Wed Jan 17 09:00:12 UTC 2024 with GitHub Actions
#
# This script refreshes the "Price" (column D) and "Volume(1h)" (column E)
# values for the crypto tickers listed on the active worksheet, matching
# the latest data pulled by the scraping job. Each price / percentage cell
# is forced to Text format ("@") before the value is written so that values
# such as "1.995.21" (dotted thousands grouping) or "10.10" / "0.999"
# (values that look numeric) are preserved exactly as plain text, instead
# of being re-interpreted as numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.725.65'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.56%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.545.71'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.19%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.47'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.31%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.56'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.14%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.50%  '

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.26%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.53'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.56%  '

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.02%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.39'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.25%  '

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.94%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.933.62'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.22%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.82'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +4.23%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.501.84'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.03%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.836'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.71%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.723.56'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.62%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.74'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.09%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.37'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.56%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0954'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.30%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.29'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.63%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '247.63'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.22%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.90'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.99%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.05'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.33%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.59'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.38%  '

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.06%  '

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.73%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.27'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.63%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.10'

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '157.75'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.01%  '

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.36%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0795'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.51%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.28'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.64%  '

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.41%  '

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -4.01%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.60'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.35%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.37'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.71%  '

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.85%  '

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.73%  '

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.03%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.07'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +5.86%  '

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.19%  '

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.61%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.995.21'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.29%  '

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.16%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.07'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.24%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.788.29'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.33%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '81.14'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.99%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.193'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.04%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.30'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.84%  '
